$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Existing rows 2-6: only touch the cells whose value actually changes
# (leave untouched cells exactly as they were).
# ------------------------------------------------------------------
$ws.Range("B2").Value = "NSE:AARTIPHARM"
$ws.Range("C2").Value = "NSE:LICMFGOLD"
$ws.Range("D2").Value = "NSE:DRREDDY"
$ws.Range("E2").Value = "NSE:HINDUNILVR"
$ws.Range("F2").Value = "NSE:NATIONALUM"

$ws.Range("B3").Value = "NSE:AARVI"
$ws.Range("D3").Value = "NSE:NATIONALUM"

$ws.Range("B4").Value = "NSE:ACI"

$ws.Range("B5").Value = "NSE:AHL"

$ws.Range("B6").Value = "NSE:AKSHARCHEM"

# ------------------------------------------------------------------
# New rows 7-30: append ticker rows, matching the formatting used by
# the existing rows (style index carried by column A: bold/centered/
# bordered "Normal" header style copied from A2).
# ------------------------------------------------------------------
$ws.Range("A2").Copy()
$ws.Range("A7:A30").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$colA = New-Object "object[,]" 24,1
$colA[0,0] = 5
$colA[1,0] = 6
$colA[2,0] = 7
$colA[3,0] = 8
$colA[4,0] = 9
$colA[5,0] = 10
$colA[6,0] = 11
$colA[7,0] = 12
$colA[8,0] = 13
$colA[9,0] = 14
$colA[10,0] = 15
$colA[11,0] = 16
$colA[12,0] = 17
$colA[13,0] = 18
$colA[14,0] = 19
$colA[15,0] = 20
$colA[16,0] = 21
$colA[17,0] = 22
$colA[18,0] = 23
$colA[19,0] = 24
$colA[20,0] = 25
$colA[21,0] = 26
$colA[22,0] = 27
$colA[23,0] = 28
$ws.Range("A7:A30").Value = $colA

$colB = New-Object "object[,]" 24,1
$colB[0,0] = "NSE:BBOX"
$colB[1,0] = "NSE:BHARATFORG"
$colB[2,0] = "NSE:BLBLIMITED"
$colB[3,0] = "NSE:CREATIVE"
$colB[4,0] = "NSE:CUPID"
$colB[5,0] = "NSE:DEEPAKNTR"
$colB[6,0] = "NSE:ECLERX"
$colB[7,0] = "NSE:EXCELINDUS"
$colB[8,0] = "NSE:GILLANDERS"
$colB[9,0] = "NSE:GREENLAM"
$colB[10,0] = "NSE:HESTERBIO"
$colB[11,0] = "NSE:INDRAMEDCO"
$colB[12,0] = "NSE:IPL"
$colB[13,0] = "NSE:JSL"
$colB[14,0] = "NSE:KPIGREEN"
$colB[15,0] = "NSE:LYKALABS"
$colB[16,0] = "NSE:NATIONALUM"
$colB[17,0] = "NSE:NMDC"
$colB[18,0] = "NSE:OMINFRAL"
$colB[19,0] = "NSE:PAKKA"
$colB[20,0] = "NSE:PDMJEPAPER"
$colB[21,0] = "NSE:PITTIENG"
$colB[22,0] = "NSE:POLYCAB"
$colB[23,0] = "NSE:PRIVISCL"
$ws.Range("B7:B30").Value = $colB

